$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the old row 5 (TRANSACTION DATE / REFERENCE / ENTRIES) down to row 6
# so we can insert the new "CLAIM ID" row at row 4 without losing data.
$ws.Range("A5:C5").Cut($ws.Range("A6:C6"))

# Clear out the now-vacated row 5 entirely (contents and formatting) so it
# stays a blank, unstyled row.
$ws.Range("A5:C5").Clear()

# Populate the new row 4 with "CLAIM ID", using the same bold style as the
# other header cells in column A (copy formatting from A3).
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "CLAIM ID"
